# Edit script: expand "groom" sheet with groom parameters table,
# and update "studio" sheet tool_state/view_state values.
# (Mirrors the authored diff: "Flipped FFC optimizeTest for better correspondence")

$wb = $excel.ActiveWorkbook

function Set-TextValue($cell, [string]$val) {
    # Force literal text storage even for values that look like booleans/numbers
    # (mirrors typing a leading apostrophe in Excel), then strip the resulting
    # quote-prefix style back off so the cell looks like a normal text cell.
    $cell.Value = "'" + $val
    $cell.Style = "Normal"
}

$groom = $wb.Worksheets.Item("groom")

$groom.Range("A2").Value = "alignment_enabled"
Set-TextValue $groom.Range("B2") "true"
Set-TextValue $groom.Range("C2") "true"

$groom.Range("A3").Value = "alignment_method"
$groom.Range("B3").Value = "Center"
$groom.Range("C3").Value = "Center"

$groom.Range("A4").Value = "antialias"
Set-TextValue $groom.Range("B4") "true"
Set-TextValue $groom.Range("C4") "true"

$groom.Range("A5").Value = "pad"
Set-TextValue $groom.Range("B5") "true"
Set-TextValue $groom.Range("C5") "true"

$groom.Range("A6").Value = "pad_value"
Set-TextValue $groom.Range("B6") "10"
Set-TextValue $groom.Range("C6") "10"

$groom.Range("A7").Value = "fastmarching"
Set-TextValue $groom.Range("B7") "true"
Set-TextValue $groom.Range("C7") "true"

$groom.Range("A8").Value = "blur"
Set-TextValue $groom.Range("B8") "true"
Set-TextValue $groom.Range("C8") "true"

$groom.Range("A9").Value = "blur_sigma"
Set-TextValue $groom.Range("B9") "2.000000"
Set-TextValue $groom.Range("C9") "2.000000"

$groom.Range("A10").Value = "isolate"
Set-TextValue $groom.Range("B10") "true"
Set-TextValue $groom.Range("C10") "true"

$groom.Range("A11").Value = "fill_holes"
Set-TextValue $groom.Range("B11") "true"
Set-TextValue $groom.Range("C11") "true"

$groom.Range("A12").Value = "fill_mesh_holes"
Set-TextValue $groom.Range("B12") "false"
Set-TextValue $groom.Range("C12") "false"

$groom.Range("A13").Value = "antialias_amount"
Set-TextValue $groom.Range("B13") "10"
Set-TextValue $groom.Range("C13") "10"

$groom.Range("A14").Value = "groom_output_prefix"
$groom.Range("B14").Value = "groomed1"
$groom.Range("C14").Value = "groomed1"

$groom.Range("A15").Value = "convert_to_mesh"
Set-TextValue $groom.Range("B15") "false"
Set-TextValue $groom.Range("C15") "false"

$groom.Range("A16").Value = "mesh_smooth"
Set-TextValue $groom.Range("B16") "false"
Set-TextValue $groom.Range("C16") "false"

$groom.Range("A17").Value = "mesh_smoothing_method"
$groom.Range("B17").Value = "Laplacian"
$groom.Range("C17").Value = "Laplacian"

$groom.Range("A18").Value = "mesh_smoothing_vtk_laplacian_iterations"
Set-TextValue $groom.Range("B18") "10"
Set-TextValue $groom.Range("C18") "10"

$groom.Range("A19").Value = "mesh_smoothing_vtk_laplacian_relaxation"
Set-TextValue $groom.Range("B19") "1.000000"
Set-TextValue $groom.Range("C19") "1.000000"

$groom.Range("A20").Value = "mesh_smoothing_vtk_windowed_sinc_iterations"
Set-TextValue $groom.Range("B20") "10"
Set-TextValue $groom.Range("C20") "10"

$groom.Range("A21").Value = "mesh_smoothing_vtk_windowed_sinc_passband"
Set-TextValue $groom.Range("B21") "0.050000"
Set-TextValue $groom.Range("C21") "0.050000"

$groom.Range("A22").Value = "crop"
Set-TextValue $groom.Range("B22") "true"
Set-TextValue $groom.Range("C22") "true"

$groom.Range("A23").Value = "reflect"
Set-TextValue $groom.Range("B23") "false"
Set-TextValue $groom.Range("C23") "false"

$groom.Range("A24").Value = "reflect_column"
$groom.Range("B24").Value = "name"
$groom.Range("C24").Value = "name"

$groom.Range("A25").Value = "reflect_choice"
$groom.Range("B25").Value = "sphere10_DT"
$groom.Range("C25").Value = "sphere10_DT"

$groom.Range("A26").Value = "reflect_axis"
$groom.Range("B26").Value = "X"
$groom.Range("C26").Value = "X"

$groom.Range("A27").Value = "resample"
Set-TextValue $groom.Range("B27") "true"
Set-TextValue $groom.Range("C27") "true"

$groom.Range("A28").Value = "isotropic"
Set-TextValue $groom.Range("B28") "true"
Set-TextValue $groom.Range("C28") "true"

$groom.Range("A29").Value = "iso_spacing"
Set-TextValue $groom.Range("B29") "1.000000"
Set-TextValue $groom.Range("C29") "1.000000"

$groom.Range("A30").Value = "spacing"
$groom.Range("B30").Value = "1 1 1"
$groom.Range("C30").Value = "1 1 1"

$groom.Range("A31").Value = "remesh"
Set-TextValue $groom.Range("B31") "true"
Set-TextValue $groom.Range("C31") "true"

$groom.Range("A32").Value = "remesh_percent_mode"
Set-TextValue $groom.Range("B32") "true"
Set-TextValue $groom.Range("C32") "true"

$groom.Range("A33").Value = "remesh_percent"
Set-TextValue $groom.Range("B33") "75.000000"
Set-TextValue $groom.Range("C33") "75.000000"

$groom.Range("A34").Value = "remesh_num_vertices"
Set-TextValue $groom.Range("B34") "3000"
Set-TextValue $groom.Range("C34") "3000"

$groom.Range("A35").Value = "remesh_gradation"
Set-TextValue $groom.Range("B35") "1.000000"
Set-TextValue $groom.Range("C35") "1.000000"

$groom.Range("A36").Value = "skip_grooming"
Set-TextValue $groom.Range("B36") "false"
Set-TextValue $groom.Range("C36") "false"

$groom.Range("A37").Value = "groom_all_domains_the_same"
Set-TextValue $groom.Range("B37") "true"
Set-TextValue $groom.Range("C37") "true"

$studio = $wb.Worksheets.Item("studio")
$studio.Range("B4").Value = "analysis"
$studio.Range("B5").Value = "Reconstructed"
